$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6493911147117615
$ws.Range("B1").Value = 0.9124007821083069
$ws.Range("C1").Value = 1.171574354171753
$ws.Range("D1").Value = 3.834018230438232
$ws.Range("E1").Value = 2.431010723114014
